$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "工业富联"
$ws.Range("B2").Value = "上海建工"
$ws.Range("C2").Value = "凯美特气"
$ws.Range("A3").Value = "立讯精密"
$ws.Range("B3").Value = "立讯精密"
$ws.Range("C3").Value = "山子高科"
$ws.Range("B4").Value = "山子高科"
$ws.Range("C4").Value = "立讯精密"
$ws.Range("A5").Value = "上海建工"
$ws.Range("B5").Value = "首开股份"
$ws.Range("C5").Value = "张江高科"
$ws.Range("A6").Value = "首开股份"
$ws.Range("B6").Value = "和而泰"
$ws.Range("A7").Value = "和而泰"
$ws.Range("B7").Value = "工业富联"
$ws.Range("A8").Value = "三花智控"
$ws.Range("B8").Value = "东华软件"
$ws.Range("C8").Value = "卧龙电驱"
$ws.Range("A9").Value = "卧龙电驱"
$ws.Range("B9").Value = "大众公用"
$ws.Range("C9").Value = "三花智控"
$ws.Range("A10").Value = "东华软件"
$ws.Range("B10").Value = "福龙马"
$ws.Range("C10").Value = "海立股份"
$ws.Range("A11").Value = "沃尔核材"
$ws.Range("B11").Value = "沃尔核材"
$ws.Range("C11").Value = "东华软件"
$ws.Range("A12").Value = "凯美特气"
$ws.Range("B12").Value = "三花智控"
$ws.Range("C12").Value = "先导智能"
$ws.Range("A13").Value = "天普股份"
$ws.Range("B13").Value = "凯美特气"
$ws.Range("C13").Value = "数据港"
$ws.Range("A14").Value = "先导智能"
$ws.Range("C14").Value = "杭电股份"
$ws.Range("B15").Value = "先导智能"
$ws.Range("C15").Value = "工业富联"
$ws.Range("A16").Value = "长川科技"
$ws.Range("B16").Value = "大洋电机"
$ws.Range("C16").Value = "波长光电"
$ws.Range("A17").Value = "英维克"
$ws.Range("B17").Value = "长川科技"
$ws.Range("C17").Value = "沃尔核材"
$ws.Range("A18").Value = "数据港"
$ws.Range("B18").Value = "万向钱潮"
$ws.Range("C18").Value = "利欧股份"
$ws.Range("A19").Value = "胜宏科技"
$ws.Range("B19").Value = "杭电股份"
$ws.Range("C19").Value = "川润股份"
$ws.Range("A20").Value = "杭电股份"
$ws.Range("B20").Value = "天普股份"
$ws.Range("C20").Value = "领益智造"
$ws.Range("A21").Value = "云南旅游"
$ws.Range("B21").Value = "东方财富"
$ws.Range("C21").Value = "天赐材料"
